# Status report update: add a new row for shipping motor/encoder RMA,
# reusing the date-formatted style from the row above (row 55) so that
# no new cell style is introduced, then move the active selection down
# to the next empty row (A57), matching how Excel behaves after a user
# types a new row of data and presses Enter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (incl. the date number format) from the last existing
# data row down into the new row 56.
$ws.Range("A55").Copy()
$ws.Range("A56").PasteSpecial(-4122)

# Fill in the new entry's values.
$ws.Cells.Item(56, 1).Value = 40245
$ws.Cells.Item(56, 2).Value = 1
$ws.Cells.Item(56, 3).Value = "RMA/ship motor/encoder"

# Move the selection to the next blank row, as Excel would after entry.
[void]$ws.Range("A57").Select()
